$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'60.331.92"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.86%  "
$ws.Range("D3").Value = "'2.346.09"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -3.53%  "
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").Value = "'540.44"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.21%  "
$ws.Range("D6").Value = "'136.33"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -5.62%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").Value = "  -10.47%  "
$ws.Range("D9").Value = "'2.347.35"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.24%  "
$ws.Range("D10").Value = "'0.104"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.15%  "
$ws.Range("E11").Value = "  +0.19%  "
$ws.Range("D12").Value = "'5.25"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.86%  "
$ws.Range("E13").Value = "  -2.38%  "
$ws.Range("D14").Value = "'24.33"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.98%  "
$ws.Range("D15").Value = "'2.769.12"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.37%  "
$ws.Range("D16").Value = "'60.345.90"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.68%  "
$ws.Range("E17").Value = "  -2.31%  "
$ws.Range("D18").Value = "'2.346.57"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -5.57%  "
$ws.Range("E19").Value = "  -3.40%  "
$ws.Range("E20").Value = "  -1.92%  "
$ws.Range("D21").Value = "'310.06"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.83%  "
$ws.Range("D22").Value = "'6.55"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -5.01%  "
$ws.Range("E23").Value = "  -0.10%  "
$ws.Range("E24").Value = "  +1.12%  "
$ws.Range("D25").Value = "'62.83"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.31%  "
$ws.Range("E26").Value = "  +11.32%  "
$ws.Range("E27").Value = "  +0.15%  "
$ws.Range("D28").Value = "'2.466.43"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.31%  "
$ws.Range("D29").Value = "'7.89"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.63%  "
$ws.Range("D30").Value = "'0.0₃0875"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -8.13%  "
$ws.Range("E31").Value = "  -4.82%  "
$ws.Range("D32").Value = "'493.01"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -5.04%  "
$ws.Range("E33").Value = "  -1.46%  "
$ws.Range("E34").Value = "  -4.57%  "
$ws.Range("E35").Value = "  -2.67%  "
$ws.Range("E36").Value = "  +0.48%  "
$ws.Range("D37").Value = "'4.54"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.87%  "
$ws.Range("E38").Value = "  -1.18%  "
$ws.Range("D39").Value = "'18.24"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.69%  "
$ws.Range("D40").Value = "'5.19"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -7.70%  "
$ws.Range("D41").Value = "'1.77"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.16%  "
$ws.Range("D43").Value = "'136.92"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.56%  "
$ws.Range("D44").Value = "'39.91"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.27%  "
$ws.Range("D45").Value = "'141.11"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.40%  "
$ws.Range("E46").Value = "  -7.09%  "
$ws.Range("D47").Value = "'3.51"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.82%  "
$ws.Range("D48").Value = "'0.0506"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.09%  "
$ws.Range("D49").Value = "'19.30"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -6.92%  "
$ws.Range("D50").Value = "'0.567"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.03%  "
$ws.Range("D51").Value = "'0.0892"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.65%  "
